$wb = $excel.ActiveWorkbook
$stations = $wb.Worksheets.Item("Stations")

# Swap lat/lng values: D should hold lng (negative), E should hold lat (positive)
for ($r = 2; $r -le 4; $r++) {
    $latVal = $stations.Cells.Item($r, 4).Value2
    $lngVal = $stations.Cells.Item($r, 5).Value2
    $stations.Cells.Item($r, 4).Value = $lngVal
    $stations.Cells.Item($r, 5).Value = $latVal
}

# Make Stations the active/selected sheet, with selection at H7
$stations.Activate() | Out-Null
$stations.Range("H7").Select() | Out-Null

